$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Update the F column timestamps (F2:F110) in the "data" sheet ---
$timeValues = @(
    "2021-10-05 14:33:17.609036",
    "2021-10-05 14:33:17.609044",
    "2021-10-05 14:33:17.609048",
    "2021-10-05 14:33:17.609050",
    "2021-10-05 14:33:17.609053",
    "2021-10-05 14:33:17.609056",
    "2021-10-05 14:33:17.609058",
    "2021-10-05 14:33:17.609061",
    "2021-10-05 14:33:17.609063",
    "2021-10-05 14:33:17.609066",
    "2021-10-05 14:33:17.609068",
    "2021-10-05 14:33:17.609071",
    "2021-10-05 14:33:17.609074",
    "2021-10-05 14:33:17.609076",
    "2021-10-05 14:33:17.609079",
    "2021-10-05 14:33:17.609081",
    "2021-10-05 14:33:17.609084",
    "2021-10-05 14:33:17.609087",
    "2021-10-05 14:33:17.609089",
    "2021-10-05 14:33:17.609092",
    "2021-10-05 14:33:17.609095",
    "2021-10-05 14:33:17.609097",
    "2021-10-05 14:33:17.609100",
    "2021-10-05 14:33:17.609102",
    "2021-10-05 14:33:17.609105",
    "2021-10-05 14:33:17.609108",
    "2021-10-05 14:33:17.609111",
    "2021-10-05 14:33:17.609113",
    "2021-10-05 14:33:17.609116",
    "2021-10-05 14:33:17.609118",
    "2021-10-05 14:33:17.609121",
    "2021-10-05 14:33:17.609123",
    "2021-10-05 14:33:17.609126",
    "2021-10-05 14:33:17.609129",
    "2021-10-05 14:33:17.609131",
    "2021-10-05 14:33:17.609134",
    "2021-10-05 14:33:17.609136",
    "2021-10-05 14:33:17.609139",
    "2021-10-05 14:33:17.609142",
    "2021-10-05 14:33:17.609144",
    "2021-10-05 14:33:17.609147",
    "2021-10-05 14:33:17.609150",
    "2021-10-05 14:33:17.609152",
    "2021-10-05 14:33:17.609155",
    "2021-10-05 14:33:17.609157",
    "2021-10-05 14:33:17.609159",
    "2021-10-05 14:33:17.609162",
    "2021-10-05 14:33:17.609164",
    "2021-10-05 14:33:17.609167",
    "2021-10-05 14:33:17.609169",
    "2021-10-05 14:33:17.609172",
    "2021-10-05 14:33:17.609174",
    "2021-10-05 14:33:17.609177",
    "2021-10-05 14:33:17.609180",
    "2021-10-05 14:33:17.609182",
    "2021-10-05 14:33:17.609185",
    "2021-10-05 14:33:17.609188",
    "2021-10-05 14:33:17.609190",
    "2021-10-05 14:33:17.609193",
    "2021-10-05 14:33:17.609195",
    "2021-10-05 14:33:17.609198",
    "2021-10-05 14:33:17.609201",
    "2021-10-05 14:33:17.609204",
    "2021-10-05 14:33:17.609206",
    "2021-10-05 14:33:17.609210",
    "2021-10-05 14:33:17.609213",
    "2021-10-05 14:33:17.609216",
    "2021-10-05 14:33:17.609219",
    "2021-10-05 14:33:17.609222",
    "2021-10-05 14:33:17.609224",
    "2021-10-05 14:33:17.609227",
    "2021-10-05 14:33:17.609230",
    "2021-10-05 14:33:17.609232",
    "2021-10-05 14:33:17.609235",
    "2021-10-05 14:33:17.609238",
    "2021-10-05 14:33:17.609240",
    "2021-10-05 14:33:17.609245",
    "2021-10-05 14:33:17.609248",
    "2021-10-05 14:33:17.609251",
    "2021-10-05 14:33:17.609253",
    "2021-10-05 14:33:17.609256",
    "2021-10-05 14:33:17.609258",
    "2021-10-05 14:33:17.609265",
    "2021-10-05 14:33:17.609269",
    "2021-10-05 14:33:17.609272",
    "2021-10-05 14:33:17.609274",
    "2021-10-05 14:33:17.609277",
    "2021-10-05 14:33:17.609279",
    "2021-10-05 14:33:17.609282",
    "2021-10-05 14:33:17.609285",
    "2021-10-05 14:33:17.609287",
    "2021-10-05 14:33:17.609290",
    "2021-10-05 14:33:17.609293",
    "2021-10-05 14:33:17.609296",
    "2021-10-05 14:33:17.609299",
    "2021-10-05 14:33:17.609302",
    "2021-10-05 14:33:17.609305",
    "2021-10-05 14:33:17.609307",
    "2021-10-05 14:33:17.609310",
    "2021-10-05 14:33:17.609313",
    "2021-10-05 14:33:17.609315",
    "2021-10-05 14:33:17.609318",
    "2021-10-05 14:33:17.609320",
    "2021-10-05 14:33:17.609323",
    "2021-10-05 14:33:17.609326",
    "2021-10-05 14:33:17.609329",
    "2021-10-05 14:33:17.609331",
    "2021-10-05 14:33:17.609334",
    "2021-10-05 14:33:17.609339"
)

for ($i = 0; $i -lt $timeValues.Length; $i++) {
    $row = $i + 2
    $dataSheet.Range("F$row").Value = $timeValues[$i]
}

# --- 2. Add the new "metadata" worksheet, placed right after "data" ---
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$meta.Name = "metadata"

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Bone Marrow Failure"
$meta.Range("C2").Value = 56
# "1.7" must stay a text value (not get coerced to the number 1.7), while
# leaving the cell on the workbook's default (unstyled) format.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.7"
$meta.Range("D2").Style = "Normal"
$meta.Range("E2").Value = "2021-09-27T07:44:31.137803Z"
$meta.Range("F2").Value = "2021-10-05 14:33:17.605641"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/56/?format=json"

# --- 3. Match the header/first-column styling used on the "data" sheet ---
$dataSheet.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

$dataSheet.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$excel.CutCopyMode = 0
[void]$meta.Range("A1").Select()

Write-Host "metadata sheet added and timestamps refreshed"
